# This workbook tracks crafting-leve profitability per job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Columns H:N hold market-price/profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# that are refreshed from current market data by a scheduled runner. This script pushes the
# refreshed figures from that run into the corresponding cells on each sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 97: Materia Worth
$ws.Range("H97").Value = 1424.8334
$ws.Range("J97").Value = 1463.4546
$ws.Range("L97").Value = 4390.3638
$ws.Range("N97").Value = -5382.3638

# Row 121: Mindful Medicine
$ws.Range("H121").Value = 2656.6
$ws.Range("I121").Value = 450
$ws.Range("J121").Value = 4127.6665
$ws.Range("K121").Value = 1350
$ws.Range("L121").Value = 12382.9995
$ws.Range("M121").Value = 397
$ws.Range("N121").Value = -15876.9995

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1963664.6
$ws.Range("I132").Value = 3601.5386
$ws.Range("J132").Value = 4087066.2
$ws.Range("K132").Value = 10804.6158
$ws.Range("L132").Value = 12261198.6
$ws.Range("M132").Value = -8274.6158
$ws.Range("N132").Value = -12266258.6

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1755665.2
$ws.Range("I137").Value = 2858263.5
$ws.Range("J137").Value = 1531.7273
$ws.Range("K137").Value = 8574790.5
$ws.Range("L137").Value = 4595.1819
$ws.Range("M137").Value = -8572240.5
$ws.Range("N137").Value = -9695.1819

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 5116.4614
$ws.Range("I141").Value = 1639.8572
$ws.Range("J141").Value = 9172.5
$ws.Range("K141").Value = 4919.571599999999
$ws.Range("L141").Value = 27517.5
$ws.Range("M141").Value = 260.4284000000007
$ws.Range("N141").Value = -37877.5

$ws = $wb.Worksheets.Item("ARM")
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 59195.055
$ws.Range("I132").Value = 46417.547
$ws.Range("J132").Value = 79274
$ws.Range("K132").Value = 139252.641
$ws.Range("L132").Value = 237822
$ws.Range("M132").Value = -136722.641
$ws.Range("N132").Value = -242882

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 989.3461
$ws.Range("I99").Value = 874.25
$ws.Range("J99").Value = 1173.5
$ws.Range("K99").Value = 874.25
$ws.Range("L99").Value = 1173.5
$ws.Range("M99").Value = 623.75
$ws.Range("N99").Value = -4169.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 3282.2
$ws.Range("I16").Value = 3005.5
$ws.Range("J16").Value = 3466.6667
$ws.Range("K16").Value = 3005.5
$ws.Range("L16").Value = 3466.6667
$ws.Range("M16").Value = -2718.5
$ws.Range("N16").Value = -4040.6667

# Row 31: Wall Not Found
$ws.Range("H31").Value = 3434.577
$ws.Range("I31").Value = 2776.6
$ws.Range("J31").Value = 4331.8184
$ws.Range("K31").Value = 2776.6
$ws.Range("L31").Value = 4331.8184
$ws.Range("M31").Value = -2481.6
$ws.Range("N31").Value = -4921.8184

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3434.577
$ws.Range("I34").Value = 2776.6
$ws.Range("J34").Value = 4331.8184
$ws.Range("K34").Value = 2776.6
$ws.Range("L34").Value = 4331.8184
$ws.Range("M34").Value = -2574.6
$ws.Range("N34").Value = -4735.8184

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 3195.5454
$ws.Range("I62").Value = 2949.6
$ws.Range("J62").Value = 3400.5
$ws.Range("K62").Value = 2949.6
$ws.Range("L62").Value = 3400.5
$ws.Range("M62").Value = -2325.6
$ws.Range("N62").Value = -4648.5

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 3195.5454
$ws.Range("I65").Value = 2949.6
$ws.Range("J65").Value = 3400.5
$ws.Range("K65").Value = 14748
$ws.Range("L65").Value = 17002.5
$ws.Range("M65").Value = -11628
$ws.Range("N65").Value = -23242.5

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1133.5834
$ws.Range("I105").Value = 860
$ws.Range("J105").Value = 1407.1666
$ws.Range("K105").Value = 860
$ws.Range("L105").Value = 1407.1666
$ws.Range("M105").Value = 887
$ws.Range("N105").Value = -4901.1666

# Row 113: Patient Patients
$ws.Range("H113").Value = 3282.2
$ws.Range("I113").Value = 3005.5
$ws.Range("J113").Value = 3466.6667
$ws.Range("K113").Value = 3005.5
$ws.Range("L113").Value = 3466.6667
$ws.Range("M113").Value = -835.5
$ws.Range("N113").Value = -7806.6667

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 32356.213
$ws.Range("I132").Value = 1806.0385
$ws.Range("J132").Value = 145828.28
$ws.Range("K132").Value = 5418.1155
$ws.Range("L132").Value = 437484.84
$ws.Range("M132").Value = -2888.1155
$ws.Range("N132").Value = -442544.84

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Range("H4").Value = 4547816
$ws.Range("I4").Value = 216.66667
$ws.Range("J4").Value = 5265858
$ws.Range("K4").Value = 650.00001
$ws.Range("L4").Value = 15797574
$ws.Range("M4").Value = -538.00001
$ws.Range("N4").Value = -15797798

# Row 12: Butter Me Up
$ws.Range("H12").Value = 31250066
$ws.Range("I12").Value = 76923144
$ws.Range("J12").Value = 63.789474
$ws.Range("K12").Value = 230769432
$ws.Range("L12").Value = 191.368422
$ws.Range("M12").Value = -230769259
$ws.Range("N12").Value = -537.368422

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 970.0513
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 977.1579
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 2931.4737
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -13011.4737

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Range("H132").Value = 102767.8
$ws.Range("I132").Value = 85788.164
$ws.Range("J132").Value = 128237.25
$ws.Range("K132").Value = 257364.492
$ws.Range("L132").Value = 384711.75
$ws.Range("M132").Value = -254834.492
$ws.Range("N132").Value = -389771.75

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 11366770
$ws.Range("I7").Value = 19232744
$ws.Range("J7").Value = 4808.6665
$ws.Range("K7").Value = 19232744
$ws.Range("L7").Value = 4808.6665
$ws.Range("M7").Value = -19232632
$ws.Range("N7").Value = -5032.6665

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 1483.1666
$ws.Range("I68").Value = 1414.9231
$ws.Range("J68").Value = 1660.6
$ws.Range("K68").Value = 1414.9231
$ws.Range("L68").Value = 1660.6
$ws.Range("M68").Value = -665.9231
$ws.Range("N68").Value = -3158.6

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 1483.1666
$ws.Range("I71").Value = 1414.9231
$ws.Range("J71").Value = 1660.6
$ws.Range("K71").Value = 7074.6155
$ws.Range("L71").Value = 8303
$ws.Range("M71").Value = -3330.6155
$ws.Range("N71").Value = -15791

# Row 126: Battered Books
$ws.Range("H126").Value = 11366770
$ws.Range("I126").Value = 19232744
$ws.Range("J126").Value = 4808.6665
$ws.Range("K126").Value = 57698232
$ws.Range("L126").Value = 14425.9995
$ws.Range("M126").Value = -57695762
$ws.Range("N126").Value = -19365.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 39: By the Short Hairs
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Row 96: Skills on Display
$ws.Range("H96").Value = 1011.625
$ws.Range("I96").Value = 925.75
$ws.Range("J96").Value = 1097.5
$ws.Range("K96").Value = 925.75
$ws.Range("L96").Value = 1097.5
$ws.Range("M96").Value = 447.25
$ws.Range("N96").Value = -3843.5

# Row 100: Of Great Import
$ws.Range("H100").Value = 67909.07000000001
$ws.Range("I100").Value = 84006.336
$ws.Range("J100").Value = 57177.555
$ws.Range("K100").Value = 168012.672
$ws.Range("L100").Value = 114355.11
$ws.Range("M100").Value = -167471.672
$ws.Range("N100").Value = -115437.11

# Row 113: A Tender Table
$ws.Range("H113").Value = 944.2
$ws.Range("I113").Value = 430.25
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1290.75
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = 879.25
$ws.Range("N113").Value = -13340
